$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 16.02.2022 12:15"

# Swap B3 (current price) and C3 (old price) values
$ws.Range("B3").Value = 35.9
$ws.Range("C3").Value = 36.5

# Update the delta sign/text in D3 - keep it stored as text (e.g. "-0.6"),
# matching the original cell which held a text value like "+0.6".
# Temporarily force a text number format so the assigned string isn't
# auto-coerced into a number, then restore the original style so no new
# style index is introduced.
$d3Style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "-0.6"
$ws.Range("D3").Style = $d3Style

# Update the timestamp text in E3
$ws.Range("E3").Value = "2022-02-16 12:15:15"
